# ---------------------------------------------------------------------------
# TC01_Login test-case workbook build-out.
#   Sheet1 "Sheet1"  -> renamed "TC01_Login", login form (UserName/Password)
#   +new sheet "0"                       -> country-selection form w/ dropdown
#   +new sheet "Resources"               -> lookup data for the dropdown list
#   + "listcountry" defined name -> Resources!$A$2:$A$3
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename existing sheet, this will become "TC01_Login" ---------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "TC01_Login"

# --- Add sheet "0" right after TC01_Login -----------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "0"

# --- Add sheet "Resources" right after "0" ----------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Resources"

# --- Defined name used by the dropdown list ---------------------------------
$wb.Names.Add('listcountry', '=Resources!$A$2:$A$3')

# =============================================================================
# Sheet "TC01_Login": UserName / Password login form
# =============================================================================
$ws1.Range("A1").Value = "UserName"
$ws1.Range("B1").Value = "user300@gmail.com"
$ws1.Range("A2").Value = "Password"
$ws1.Range("B2").Value = "123456789oO"

$ws1.Range("A1:B2").VerticalAlignment = -4108

$ws1.Columns.Item(1).ColumnWidth = 14.833333333333334
$ws1.Columns.Item(2).ColumnWidth = 31.333333333333332

$ws1.Rows.Item(1).RowHeight = 48.75
$ws1.Rows.Item(2).RowHeight = 55.5

$ws1.Range("B2").Select()

# =============================================================================
# Sheet "0": country-name field with dropdown validation
# =============================================================================
$ws2.Range("A1").Value = "Country Name "
$ws2.Range("B1").Value = "France"

$ws2.Columns.Item(1).ColumnWidth = 26.666666666666668
$ws2.Columns.Item(2).ColumnWidth = 53.166666666666664

$ws2.Rows.Item(1).RowHeight = 48
$ws2.Rows.Item(2).RowHeight = 48.75

$ws2.Range("B1").Validation.Add(3, 1, 1, "listcountry")

# =============================================================================
# Sheet "Resources": lookup data backing the dropdown
# =============================================================================
$ws3.Range("A1").Value = "List countries"
$ws3.Range("A2").Value = "France"
$ws3.Range("A3").Value = "United States"

$ws3.Columns.Item(1).ColumnWidth = 17.666666666666668

$ws3.Range("A2:A3").Select()

# =============================================================================
# Final active sheet/tab: sheet "0" (activeTab = 1)
# =============================================================================
$ws2.Select()
$ws2.Range("K7").Select()
